$wb = $excel.ActiveWorkbook

$oldId = "6328dcac-9d07-4ccf-a8e1-3389ef3f5daa"
$newId = "a04b0430-406c-4e47-9e35-46b35a874fe4"

$oldMd       = "$oldId.md"
$newMd       = "$newId.md"
$oldPath     = "e2e\$oldId.md"
$newPath     = "e2e\$newId.md"

$oldZhXlf    = "$oldId.264d6596c94c446420783b680b22d772fb1f7dee.zh-cn.xlf"
$newZhXlf    = "$newId.3914c112f2ce338a9c100e4757515d7effef857b.zh-cn.xlf"
$oldDeXlf    = "$oldId.264d6596c94c446420783b680b22d772fb1f7dee.de-de.xlf"
$newDeXlf    = "$newId.3914c112f2ce338a9c100e4757515d7effef857b.de-de.xlf"

$oldHoDate   = "2016-08-17 08:56:36"
$newHoDate   = "2016-08-17 08:56:52"
$oldZhDate   = "2016-08-17 08:56:31"
$newZhDate   = "2016-08-17 08:56:47"

# The hyperlink address (stored in the worksheet's external relationship) is
# left untouched by this change - only the visible display text / cell text
# gets the new generated-report file name.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/095a35a4d9e36d286634fda9bf134e62f8936864/e2e/$oldId.md"

function Update-HyperlinkDisplay($ws, $cellAddr, $newText) {
    $ws.Hyperlinks.Delete()
    $ws.Range($cellAddr).Value = $newText
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $linkAddress, "", "", $newText) | Out-Null
}

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
Update-HyperlinkDisplay $wsOverview "B2" $newPath
$wsOverview.Range("G2").Value = $newHoDate

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-HyperlinkDisplay $wsZh "A2" $newMd
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhDate

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")
Update-HyperlinkDisplay $wsDe "A2" $newMd
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHoDate
